$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "96.347.45"
$ws.Range("E2").Value = "  +0.91%  "
$ws.Range("D3").Value = "3.581.67"
$ws.Range("E3").Value = "  -0.26%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "241.28"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.18%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "654.66"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.04%  "
$ws.Range("E7").Value = "  +5.97%  "
$ws.Range("E8").Value = "  -0.39%  "
$ws.Range("E9").Value = "  +0.06%  "
$ws.Range("E10").Value = "  +3.20%  "
$ws.Range("D11").Value = "3.579.75"
$ws.Range("E11").Value = "  -0.24%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "43.13"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.20%  "
$ws.Range("E13").Value = "  +0.59%  "
$ws.Range("E14").Value = "  +1.01%  "
$ws.Range("D15").Value = "4.248.14"
$ws.Range("E15").Value = "  -0.24%  "
$ws.Range("D16").Value = "96.278.59"
$ws.Range("E16").Value = "  +1.01%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.0000259"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +1.41%  "
$ws.Range("D18").Value = "3.577.69"
$ws.Range("E18").Value = "  -0.34%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.74"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -5.43%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "12.50"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.28%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "17.72"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.92%  "
$ws.Range("E22").Value = "  +1.57%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "511.22"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.13%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "3.42"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -2.94%  "
$ws.Range("E25").Value = "  +3.75%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "6.82"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +2.66%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "96.37"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.51%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "12.63"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.78%  "
$ws.Range("D29").Value = "3.774.78"
$ws.Range("E29").Value = "  -0.17%  "
$ws.Range("B30").Value = "PancakeSwap"
$ws.Range("C30").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.98"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -7.18%  "
$ws.Range("B31").Value = "Hedera"
$ws.Range("C31").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.149"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +7.11%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "11.44"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.99%  "
$ws.Range("E33").Value = "  +0.11%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.182"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +3.44%  "
$ws.Range("E35").Value = "  +1.32%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "31.57"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -1.32%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "615.93"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +6.56%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.565"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.84%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "8.63"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +4.10%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.62"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +9.37%  "
$ws.Range("E42").Value = "  -0.09%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.906"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -2.38%  "
$ws.Range("E44").Value = "  +5.37%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "5.69"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.68%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.28"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +1.33%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "34.15"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +1.30%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "23.51"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.04%  "
$ws.Range("E49").Value = "  -0.39%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "3.60"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +3.75%  "
$ws.Range("B51").Value = "Cosmos"
$ws.Range("C51").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "8.17"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.53%  "
